$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.028.45"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "3.740.67"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'601.73"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").Value = "'167.09"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D7").Value = "3.739.78"
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("E10").Value = "  +2.44%  "
$ws.Range("E11").Value = "  +2.91%  "
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "'37.95"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("E14").Value = "  +1.86%  "
$ws.Range("D15").Value = "4.368.13"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").Value = "3.750.49"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("D17").Value = "69.012.23"
$ws.Range("D18").Value = "'7.35"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("D20").Value = "'17.27"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").Value = "'11.00"
$ws.Range("E21").Value = "  +19.51%  "
$ws.Range("D22").Value = "'492.45"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").Value = "'0.726"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("E24").Value = "  +8.29%  "
$ws.Range("D25").Value = "'84.77"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").Value = "'12.32"
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  +2.15%  "
$ws.Range("E31").Value = "  +4.33%  "
$ws.Range("E32").Value = "  +2.32%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "3.885.82"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").Value = "3.677.17"
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  +1.15%  "
$ws.Range("D39").Value = "'5.94"
$ws.Range("E39").Value = "  +3.42%  "
$ws.Range("D40").Value = "'0.135"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("D42").Value = "'2.98"
$ws.Range("E42").Value = "  +5.89%  "
$ws.Range("D43").Value = "'431.47"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").Value = "'48.70"
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "'40.31"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").Value = "'141.24"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("D50").Value = "2.784.00"
$ws.Range("E50").Value = "  +1.76%  "
$ws.Range("E51").Value = "  +0.78%  "
